$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"

$ws.Range("A7").Value = "may"
$ws.Range("B7").Value = "levi"
$ws.Range("C7").Value = "123456789"
$ws.Range("D7").Value = "qiryat gat"
$ws.Range("E7").Value = "10/10/94"
$ws.Range("F7").Value = "054-6337879"
